$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new PD1 (NAT105-CellMarque) example image paths in column D,
# rows 2-11, mirroring the existing HE/PDL1 path pattern.
$ws.Range("D2").Value  = "\PD1(NAT105)-CellMarque\PD1(NAT105)-CellMarque_A12_v3_b3_000.jpg"
$ws.Range("D3").Value  = "\PD1(NAT105)-CellMarque\PD1(NAT105)-CellMarque_A12_v3_b3_001.jpg"
$ws.Range("D4").Value  = "\PD1(NAT105)-CellMarque\PD1(NAT105)-CellMarque_A12_v3_b3_002.jpg"
$ws.Range("D5").Value  = "\PD1(NAT105)-CellMarque\PD1(NAT105)-CellMarque_A12_v3_b3_003.jpg"
$ws.Range("D6").Value  = "\PD1(NAT105)-CellMarque\PD1(NAT105)-CellMarque_A12_v3_b3_004.jpg"
$ws.Range("D7").Value  = "\PD1(NAT105)-CellMarque\PD1(NAT105)-CellMarque_A12_v3_b3_005.jpg"
$ws.Range("D8").Value  = "\PD1(NAT105)-CellMarque\PD1(NAT105)-CellMarque_A12_v3_b3_006.jpg"
$ws.Range("D9").Value  = "\PD1(NAT105)-CellMarque\PD1(NAT105)-CellMarque_A12_v3_b3_007.jpg"
$ws.Range("D10").Value = "\PD1(NAT105)-CellMarque\PD1(NAT105)-CellMarque_A12_v3_b3_008.jpg"
$ws.Range("D11").Value = "\PD1(NAT105)-CellMarque\PD1(NAT105)-CellMarque_A12_v3_b3_009.jpg"

# The header row (A1:E1) loses its bold/bordered "header" style once the new
# column gets edited - reset it back to the workbook's default "Normal" style.
$ws.Range("A1:E1").Style = "Normal"

# Column D is now much wider to fit the new (longer) PD1 path strings -
# widen it (best-fit) to fit the new content.
$ws.Columns.Item(4).ColumnWidth = 66.67

# Move the active selection to F5.
$ws.Range("F5").Select() | Out-Null
